$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 385
$ws.Cells.Item(385, 4).Value = 22

# Row 386
$ws.Cells.Item(386, 4).Value = 62
$ws.Cells.Item(386, 5).Value = 39.9
$ws.Cells.Item(386, 8).Value = 10.3

# Row 387
$ws.Cells.Item(387, 4).Value = 113
$ws.Cells.Item(387, 5).Value = 44.1
$ws.Cells.Item(387, 8).Value = 8.3

# Row 388
$ws.Cells.Item(388, 4).Value = 156
$ws.Cells.Item(388, 5).Value = 50
$ws.Cells.Item(388, 8).Value = 6.2

# Row 389
$ws.Cells.Item(389, 4).Value = 44
$ws.Cells.Item(389, 5).Value = 57.7
$ws.Cells.Item(389, 8).Value = 12.8

# Row 390
$ws.Cells.Item(390, 4).Value = 62
$ws.Cells.Item(390, 5).Value = 56.8
$ws.Cells.Item(390, 8).Value = 14.7

# Row 391
$ws.Cells.Item(391, 4).Value = 40
$ws.Cells.Item(391, 5).Value = 46.8
$ws.Cells.Item(391, 6).Value = 0.016
$ws.Cells.Item(391, 8).Value = 16.7

# Row 392
$ws.Cells.Item(392, 4).Value = 60
$ws.Cells.Item(392, 5).Value = 45.5
$ws.Cells.Item(392, 6).Value = 0
$ws.Cells.Item(392, 8).Value = 7

# Row 393
$ws.Cells.Item(393, 4).Value = 58
$ws.Cells.Item(393, 5).Value = 54.7
$ws.Cells.Item(393, 8).Value = 12.6

# Row 394
$ws.Cells.Item(394, 5).Value = 57.6
$ws.Cells.Item(394, 6).Value = 0.15
$ws.Cells.Item(394, 8).Value = 10.4

# Row 395
$ws.Cells.Item(395, 5).Value = 50.9
$ws.Cells.Item(395, 6).Value = 0.15
$ws.Cells.Item(395, 8).Value = 9.8

# Row 396
$ws.Cells.Item(396, 4).Value = 60
$ws.Cells.Item(396, 5).Value = 52.3
$ws.Cells.Item(396, 6).Value = 0.25
$ws.Cells.Item(396, 8).Value = 11.7

# Row 397
$ws.Cells.Item(397, 4).Value = 33
$ws.Cells.Item(397, 5).Value = 56.3
$ws.Cells.Item(397, 6).Value = 0.4
$ws.Cells.Item(397, 8).Value = 10.6

# Row 398
$ws.Cells.Item(398, 4).Value = 49
$ws.Cells.Item(398, 5).Value = 56.5
$ws.Cells.Item(398, 6).Value = 0
$ws.Cells.Item(398, 8).Value = 0.2

# Row 399
$ws.Cells.Item(399, 4).Value = 71
$ws.Cells.Item(399, 5).Value = 55

# Row 400
$ws.Cells.Item(400, 4).Value = 55
$ws.Cells.Item(400, 6).Value = 0.05

# Row 403
$ws.Cells.Item(403, 4).Value = 62

# Row 404
$ws.Cells.Item(404, 4).Value = 36

# Row 406
$ws.Cells.Item(406, 4).Value = 63

# Row 407
$ws.Cells.Item(407, 4).Value = 54
